$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53
$ws.Range("A53").Value = 111901547
$ws.Range("Q53").Value = 477523.7983399219
$ws.Range("R53").Value = 7033329.5960882
$ws.Range("AC53").Value = "ringhack"

# Row 54
$ws.Range("A54").Value = 111901550
$ws.Range("Q54").Value = 477473.2080285564
$ws.Range("R54").Value = 7033403.83150613

# Row 55
$ws.Range("A55").Value = 111901544
$ws.Range("Q55").Value = 477638.5281090657
$ws.Range("R55").Value = 7033514.606123095

# Row 57
$ws.Range("A57").Value = 111901585
$ws.Range("Q57").Value = 478338.6535977835
$ws.Range("R57").Value = 7035076.017275342

# Row 58
$ws.Range("A58").Value = 111901551
$ws.Range("Q58").Value = 477432.586959724
$ws.Range("R58").Value = 7033429.191801991
$ws.Range("AC58").Value = "ringhack"

# Row 59
$ws.Range("A59").Value = 111901619
$ws.Range("B59").Value = 85062
$ws.Range("E59").Value = 249278
$ws.Range("F59").Value = "Barrviolspindling"
$ws.Range("G59").Value = "Cortinarius harcynicus"
$ws.Range("H59").Value = "(Pers.) M.M.Moser"
$ws.Range("K59").Value = ""
$ws.Range("L59").Value = ""
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = ""
$ws.Range("Q59").Value = 478523.10621621
$ws.Range("R59").Value = 7034650.501669589
$ws.Range("AC59").Value = ""

# Row 60
$ws.Range("A60").Value = 111901519
$ws.Range("B60").Value = 86223
$ws.Range("E60").Value = 4412
$ws.Range("F60").Value = "Äggvaxskivling"
$ws.Range("G60").Value = "Hygrophorus karstenii"
$ws.Range("H60").Value = "Sacc. & Cub."
$ws.Range("K60").Value = ""
$ws.Range("L60").Value = ""
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = ""
$ws.Range("Q60").Value = 477765.0897337386
$ws.Range("R60").Value = 7033404.474773662
$ws.Range("AC60").Value = ""

# Row 61
$ws.Range("A61").Value = 111901545
$ws.Range("B61").Value = 56398
$ws.Range("E61").Value = 100109
$ws.Range("F61").Value = "Tretåig hackspett"
$ws.Range("G61").Value = "Picoides tridactylus"
$ws.Range("H61").Value = "(Linnaeus, 1758)"
$ws.Range("K61").Value = ""
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = ""
$ws.Range("Q61").Value = 477666.6718496145
$ws.Range("R61").Value = 7033500.081917049
$ws.Range("AC61").Value = "ringhack äldre"

# Row 62
$ws.Range("A62").Value = 111901546
$ws.Range("B62").Value = 56398
$ws.Range("E62").Value = 100109
$ws.Range("F62").Value = "Tretåig hackspett"
$ws.Range("G62").Value = "Picoides tridactylus"
$ws.Range("H62").Value = "(Linnaeus, 1758)"
$ws.Range("I62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("Q62").Value = 477668.4830064432
$ws.Range("R62").Value = 7033374.253324097
$ws.Range("AC62").Value = "ringhack äldre"

# Row 63
$ws.Range("A63").Value = 111901518
$ws.Range("B63").Value = 86223
$ws.Range("E63").Value = 4412
$ws.Range("F63").Value = "Äggvaxskivling"
$ws.Range("G63").Value = "Hygrophorus karstenii"
$ws.Range("H63").Value = "Sacc. & Cub."
$ws.Range("Q63").Value = 477673.8480424859
$ws.Range("R63").Value = 7033500.479530043

# Row 64
$ws.Range("A64").Value = 111901549
$ws.Range("Q64").Value = 477463.5153726833
$ws.Range("R64").Value = 7033364.497689161
$ws.Range("AC64").Value = "ringhack färska"

# Row 65
$ws.Range("A65").Value = 111901548
$ws.Range("B65").Value = 56398
$ws.Range("E65").Value = 100109
$ws.Range("F65").Value = "Tretåig hackspett"
$ws.Range("G65").Value = "Picoides tridactylus"
$ws.Range("H65").Value = "(Linnaeus, 1758)"
$ws.Range("K65").Value = ""
$ws.Range("L65").Value = ""
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("Q65").Value = 477476.2138649083
$ws.Range("R65").Value = 7033385.004830683
$ws.Range("AC65").Value = "ringhack äldre"

# Row 66
$ws.Range("A66").Value = 111901587
$ws.Range("B66").Value = 56543
$ws.Range("E66").Value = 103021
$ws.Range("F66").Value = "Talltita"
$ws.Range("G66").Value = "Poecile montanus"
$ws.Range("H66").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I66").NumberFormat = "@"
$ws.Range("I66").Value = "2"
$ws.Range("N66").Value = "observerad"
$ws.Range("Q66").Value = 477611.096214832
$ws.Range("R66").Value = 7033310.625431053
$ws.Range("AC66").Value = ""
